# API Fluent Entities Config Part 1
# Rework of the "Permission/Actions" tables into "Permissions/Activities",
# plus assorted FK-highlight formatting passes on the DB schema sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors as Long (RGB(r,g,b) -> r + g*256 + b*65536, matches Excel's .Color property)
function RGBColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$BlueColor = RGBColor 0 176 240    # FF00B0F0
$GreenColor = RGBColor 146 208 80  # FF92D050

# --- Rename table / header labels -------------------------------------------------
$ws.Range("G1").Value2 = "permissions"
$ws.Range("I1").Value2 = "activities"
$ws.Range("A4").Value2 = "Permissions"
$ws.Range("A5").Value2 = "Activities"

# --- New columns in existing tables ------------------------------------------------
$ws.Range("J4").Value2 = "Type"
$ws.Range("J4").Font.Color = 255          # match existing "field name" red style
$ws.Range("I5").Value2 = "FunctionId"

# --- Fix typo in categoriesTrans table ---------------------------------------------
$ws.Range("H18").Value2 = "Description"

# --- Reorder postTrans table columns (LanguageId moves to the end) -----------------
$ws.Range("I19").Value2 = "SEODescription"
$ws.Range("I20").Value2 = "SEOTitle"
$ws.Range("I21").Value2 = "SEOAlias"
$ws.Range("I22").Value2 = "LanguageId"

# --- Fix Log table column name ------------------------------------------------------
$ws.Range("D29").Value2 = "ActionId"

# --- Apply blue "FK reference" font color ---------------------------------------------
$blueCells = @("F2","G2","F3","G3","G4","I5","F16","H16","I16","E18","I22","E23","H23","D29","D31")
foreach ($addr in $blueCells) {
    $ws.Range($addr).Font.Color = $BlueColor
}

# --- Apply green "self FK reference" font color ----------------------------------------
$greenCells = @("H6","G7","D16","F19")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Font.Color = $GreenColor
}

# --- Column width adjustments --------------------------------------------------------
# Column G (7) becomes its own width group (14), same width as column E; column F (6)
# is left untouched so it keeps its original width (10.796875).
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# --- Selection matches the authored commit -------------------------------------------
$ws.Range("I11").Select()
